$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column P (new "2022" column) -----------------------------------------
# Seed every touched row's formatting from the matching column-N cell
# (same row styling as the rest of the table) before writing values, so the
# new column inherits the existing look without disturbing shared styles.

# Row 3 (bottom border separator row) - stays empty, just needs formatting.
$ws.Range("N3").Copy()
$ws.Range("P3").PasteSpecial(-4122) # xlPasteFormats

# Row 4 (year header row) - new year 2022.
$ws.Range("N4").Copy()
$ws.Range("P4").PasteSpecial(-4122)
$ws.Range("P4").Value = 2022

# Row 5 (section header row) - stays empty, just needs formatting.
$ws.Range("N5").Copy()
$ws.Range("P5").PasteSpecial(-4122)

# Row 6 (data row) - new value.
$ws.Range("N6").Copy()
$ws.Range("P6").PasteSpecial(-4122)
$ws.Range("P6").Value = 1373

# Row 7 (data row) - new value is a dash, right aligned (new cell style).
$ws.Range("N7").Copy()
$ws.Range("P7").PasteSpecial(-4122)
$ws.Range("P7").Value = "-"
$ws.Range("P7").HorizontalAlignment = -4152 # xlRight

# Row 8 (data row) - new value.
$ws.Range("N8").Copy()
$ws.Range("P8").PasteSpecial(-4122)
$ws.Range("P8").Value = 117

# Row 9 (data row) - new value.
$ws.Range("N9").Copy()
$ws.Range("P9").PasteSpecial(-4122)
$ws.Range("P9").Value = 154

# Row 10 (last data row, bottom border) - new value.
$ws.Range("N10").Copy()
$ws.Range("P10").PasteSpecial(-4122)
$ws.Range("P10").Value = 885

# Update the remembered selection to match the authored edit (P7).
[void]$ws.Range("P7").Select()
